$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Variable label used across all rows
$varLabel = "Taxa de desocupação, na semana de referência, das pessoas de 14 anos ou mais de idade"

# Final dataset (in desired row order) for rows 2..40
$data = @(
    @("Brasil","01/10/2012",6.9),
    @("Brasil","01/10/2013",6.3),
    @("Brasil","01/10/2014",6.6),
    @("Brasil","01/10/2015",9.1),
    @("Brasil","01/10/2016",12.2),
    @("Brasil","01/10/2017",11.9),
    @("Brasil","01/10/2018",11.7),
    @("Brasil","01/10/2019",11.1),
    @("Brasil","01/10/2020",14.2),
    @("Brasil","01/10/2021",11.1),
    @("Brasil","01/10/2022",7.9),
    @("Brasil","01/10/2023",7.4),
    @("Brasil","01/10/2024",6.2),
    @("Nordeste","01/10/2012",9.4),
    @("Nordeste","01/10/2013",8),
    @("Nordeste","01/10/2014",8.300000000000001),
    @("Nordeste","01/10/2015",10.6),
    @("Nordeste","01/10/2016",14.5),
    @("Nordeste","01/10/2017",14),
    @("Nordeste","01/10/2018",14.5),
    @("Nordeste","01/10/2019",13.8),
    @("Nordeste","01/10/2020",0),
    @("Nordeste","01/10/2021",0),
    @("Nordeste","01/10/2022",10.9),
    @("Nordeste","01/10/2023",10.4),
    @("Nordeste","01/10/2024",8.6),
    @("Sergipe","01/10/2012",9.699999999999999),
    @("Sergipe","01/10/2013",8.800000000000001),
    @("Sergipe","01/10/2014",9),
    @("Sergipe","01/10/2015",10.1),
    @("Sergipe","01/10/2016",15.2),
    @("Sergipe","01/10/2017",13.6),
    @("Sergipe","01/10/2018",15.2),
    @("Sergipe","01/10/2019",15),
    @("Sergipe","01/10/2020",0),
    @("Sergipe","01/10/2021",0),
    @("Sergipe","01/10/2022",11.9),
    @("Sergipe","01/10/2023",11.2),
    @("Sergipe","01/10/2024",8.4)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $varLabel
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $row++
}
